# Apply latest crypto price/volume snapshot to the sheet (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.240.83"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "'1.600.99"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'212.30"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").Value = "'18.10"
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").Value = "'0.0812"
$ws.Range("E11").Value = "  +1.87%  "
$ws.Range("D12").Value = "'1.822.21"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "'1.601.15"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").Value = "'4.04"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("D16").Value = "'26.221.31"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "'61.24"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "'202.17"
$ws.Range("E20").Value = "  +1.66%  "
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("E24").Value = "  +9.22%  "
$ws.Range("D25").Value = "'144.24"
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  -7.93%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("D30").Value = "'0.0488"
$ws.Range("E30").Value = "  +3.41%  "
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("E33").Value = "  -3.07%  "
$ws.Range("E34").Value = "  +2.69%  "
$ws.Range("E35").Value = "  -1.18%  "
$ws.Range("D36").Value = "'1.157.07"
$ws.Range("E36").Value = "  +4.42%  "
$ws.Range("E37").Value = "  +8.41%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'0.791"
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.32"
$ws.Range("E40").Value = "  -1.15%  "
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "'5.22"
$ws.Range("E43").Value = "  +1.95%  "
$ws.Range("D44").Value = "'1.737.19"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").Value = "'91.86"
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("E46").Value = "  -2.90%  "
$ws.Range("D47").Value = "'54.06"
$ws.Range("E47").Value = "  +0.96%  "
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.407"
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.0₇0955"
$ws.Range("E50").Value = "  -9.55%  "
$ws.Range("E51").Value = "  -0.16%  "
